# Update Excel SCD0011 until SCD0016
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from SCD0209 to SCD0011
$ws.Name = "SCD0011"

# Update TC_ID cell (B2): DGS-224 -> SCD0011-040
$ws.Range("B2").Value = "SCD0011-040"

# Update selection to B3
$ws.Range("B3").Select()

# Update column B width (target stored width 12.42578125; engine quantizes
# internally, 11.65 lands mid-bucket on the closest reachable stored width)
$ws.Columns.Item(2).ColumnWidth = 11.65
